# Update Betfair Back/Lay odds values on Sheet1 to reflect the latest refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Egyptian Premier: National Bank vs Zamalek
$ws.Range("F3").Value = 3.05
$ws.Range("H3").Value = 2.3
$ws.Range("I3").Value = 2.74
$ws.Range("J3").Value = 2.76
$ws.Range("K3").Value = 4.2

# Row 4 - Italian Serie A: Cagliari vs Sassuolo
$ws.Range("N4").Value = 3.2
$ws.Range("O4").Value = 1.43
$ws.Range("P4").Value = 1.71
$ws.Range("Z4").Value = 19.5
$ws.Range("AB4").Value = 10
$ws.Range("AC4").Value = 7
$ws.Range("AJ4").Value = 46
$ws.Range("AL4").Value = 55

# Row 7 - Swiss Super League: Grasshoppers Zurich vs Young Boys
$ws.Range("F7").Value = 3.35
$ws.Range("G7").Value = 3.85
$ws.Range("H7").Value = 1.98
$ws.Range("I7").Value = 2.2
$ws.Range("J7").Value = 3.8
$ws.Range("P7").Value = 2.34
$ws.Range("Q7").Value = 1.58

# Row 8 - Swiss Super League: Lugano vs Luzern
$ws.Range("H8").Value = 4.1
$ws.Range("I8").Value = 4.7
$ws.Range("K8").Value = 4.5

# Row 9 - Italian Serie A: Pisa vs Lazio
$ws.Range("F9").Value = 3.95
$ws.Range("I9").Value = 2.24
$ws.Range("K9").Value = 3.4
$ws.Range("N9").Value = 3.05
$ws.Range("O9").Value = 1.45
$ws.Range("P9").Value = 1.67
$ws.Range("AE9").Value = 28
$ws.Range("AL9").Value = 75
$ws.Range("AN9").Value = 75

# Row 12 - Colombian Primera A: Union Magdalena vs Tolima
$ws.Range("F12").Value = 3.4
$ws.Range("I12").Value = 2.32
$ws.Range("J12").Value = 3.15
$ws.Range("K12").Value = 3.6
$ws.Range("Q12").Value = 2.26

# Row 13 - Colombian Primera A: Boyaca Chico vs America de Cali S.A
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 1.87
$ws.Range("K13").Value = 3.95
$ws.Range("Q13").Value = 2.36
